$d = $word.ActiveDocument

$oldText = "Think of, or find examples from the Bible where people showed, felt or expressed love for another person or towards God."
$newText = "Read Psalm 51 – Describe the different ways that love is expressed in this Psalm.  How can you apply them to your life?"

# Locate the paragraph that currently holds both the "Think of..." text
# and the trailing _GoBack bookmark, then strip the text out of it so the
# bookmark is left behind alone in its own (now empty) paragraph.
$findRng = $d.Content
$findRng.Find.ClearFormatting()
$found = $findRng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target paragraph text"
}
$bookmarkParaIndex = $findRng.Paragraphs.Item(1).Index
$findRng.Text = ""

# Add a brand-new paragraph right after the (now bookmark-only) paragraph and
# put the "Think of..." text into it.
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)
$bookmarkPara.Range.InsertParagraphAfter()

$thinkParaIndex = $bookmarkParaIndex + 1
$thinkPara = $d.Paragraphs.Item($thinkParaIndex)
$thinkPara.Range.Text = $oldText

# Insert four blank paragraphs, then a final paragraph with the new
# "Read Psalm 51..." prompt, all immediately after the "Think of..." paragraph.
$thinkPara = $d.Paragraphs.Item($thinkParaIndex)
$insertPoint = $thinkPara.Range
$insertPoint.InsertParagraphAfter()
$insertPoint.InsertParagraphAfter()
$insertPoint.InsertParagraphAfter()
$insertPoint.InsertParagraphAfter()
$insertPoint.InsertParagraphAfter()

$psalmParaIndex = $thinkParaIndex + 5
$psalmPara = $d.Paragraphs.Item($psalmParaIndex)
$psalmPara.Range.Text = $newText
